$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 19092.691
$ws.Range("I62").Value = 5424.75
$ws.Range("K62").Value = 5424.75
$ws.Range("M62").Value = -4800.75
$ws.Range("H65").Value = 19092.691
$ws.Range("I65").Value = 5424.75
$ws.Range("K65").Value = 27123.75
$ws.Range("M65").Value = -24003.75
$ws.Range("H76").Value = 6484.9
$ws.Range("I76").Value = 4462.25
$ws.Range("K76").Value = 4462.25
$ws.Range("M76").Value = -4147.25
$ws.Range("H79").Value = 6484.9
$ws.Range("I79").Value = 4462.25
$ws.Range("K79").Value = 4462.25
$ws.Range("M79").Value = -3370.25
$ws.Range("H137").Value = 5110.6206
$ws.Range("J137").Value = 7540.636
$ws.Range("L137").Value = 22621.908
$ws.Range("N137").Value = -27721.908

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9807978
$ws.Range("I32").Value = 11365633
$ws.Range("K32").Value = 11365633
$ws.Range("M32").Value = -11365346
$ws.Range("H61").Value = 68188790
$ws.Range("I61").Value = 62506450
$ws.Range("K61").Value = 62506450
$ws.Range("M61").Value = -62506238
$ws.Range("H63").Value = 5992.6665
$ws.Range("I63").Value = 2650
$ws.Range("K63").Value = 2650
$ws.Range("M63").Value = -1964
$ws.Range("H66").Value = 5992.6665
$ws.Range("I66").Value = 2650
$ws.Range("K66").Value = 13250
$ws.Range("M66").Value = -9818
$ws.Range("H123").Value = 48710
$ws.Range("J123").Value = 48710
$ws.Range("L123").Value = 48710
$ws.Range("N123").Value = -58510
$ws.Range("H132").Value = 12422.6
$ws.Range("I132").Value = 8255.223
$ws.Range("K132").Value = 24765.669
$ws.Range("M132").Value = -22235.669
$ws.Range("H136").Value = 68188790
$ws.Range("I136").Value = 62506450
$ws.Range("K136").Value = 187519350
$ws.Range("M136").Value = -187516800

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3880.1667
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H86").Value = 4601.6665
$ws.Range("I86").Value = 4601.6665
$ws.Range("K86").Value = 4601.6665
$ws.Range("M86").Value = -3478.6665
$ws.Range("H88").Value = 29993
$ws.Range("J88").Value = 29994.5
$ws.Range("L88").Value = 29994.5
$ws.Range("N88").Value = -30806.5
$ws.Range("H89").Value = 4601.6665
$ws.Range("I89").Value = 4601.6665
$ws.Range("K89").Value = 23008.3325
$ws.Range("M89").Value = -17392.3325
$ws.Range("H91").Value = 29993
$ws.Range("J91").Value = 29994.5
$ws.Range("L91").Value = 29994.5
$ws.Range("N91").Value = -32802.5
$ws.Range("H94").Value = 1582.5883
$ws.Range("I94").Value = 1365.4
$ws.Range("K94").Value = 1365.4
$ws.Range("M94").Value = -914.4000000000001
$ws.Range("H105").Value = 3276
$ws.Range("J105").Value = 2155.3333
$ws.Range("L105").Value = 2155.3333
$ws.Range("N105").Value = -5649.3333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 35202.4
$ws.Range("J11").Value = 43753
$ws.Range("L11").Value = 43753
$ws.Range("N11").Value = -44033
$ws.Range("H12").Value = 4000
$ws.Range("J12").Value = 4000
$ws.Range("L12").Value = 4000
$ws.Range("N12").Value = -4340
$ws.Range("H31").Value = 656761.1
$ws.Range("I31").Value = 10332.263
$ws.Range("K31").Value = 10332.263
$ws.Range("M31").Value = -10037.263
$ws.Range("H34").Value = 656761.1
$ws.Range("I34").Value = 10332.263
$ws.Range("K34").Value = 10332.263
$ws.Range("M34").Value = -10130.263
$ws.Range("H62").Value = 561942
$ws.Range("I62").Value = 775359.1
$ws.Range("K62").Value = 775359.1
$ws.Range("M62").Value = -774735.1
$ws.Range("H65").Value = 561942
$ws.Range("I65").Value = 775359.1
$ws.Range("K65").Value = 3876795.5
$ws.Range("M65").Value = -3873675.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 170.61765
$ws.Range("I2").Value = 117.90909
$ws.Range("J2").Value = 195.82608
$ws.Range("K2").Value = 707.4545400000001
$ws.Range("L2").Value = 1174.95648
$ws.Range("M2").Value = -594.4545400000001
$ws.Range("N2").Value = -1400.95648
$ws.Range("H107").Value = 609.5
$ws.Range("I107").Value = 489.6875
$ws.Range("J107").Value = 849.125
$ws.Range("K107").Value = 1469.0625
$ws.Range("L107").Value = 2547.375
$ws.Range("M107").Value = 450.9375
$ws.Range("N107").Value = -6387.375
$ws.Range("H113").Value = 1166.1111
$ws.Range("I113").Value = 749.5
$ws.Range("J113").Value = 1218.1875
$ws.Range("K113").Value = 2248.5
$ws.Range("L113").Value = 3654.5625
$ws.Range("M113").Value = -78.5
$ws.Range("N113").Value = -7994.5625
$ws.Range("H131").Value = 4907.926
$ws.Range("I131").Value = 6001.4443
$ws.Range("J131").Value = 4689.222
$ws.Range("K131").Value = 18004.3329
$ws.Range("L131").Value = 14067.666
$ws.Range("M131").Value = -12964.3329
$ws.Range("N131").Value = -24147.666

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7280625.5
$ws.Range("J11").Value = 5898572.5
$ws.Range("L11").Value = 5898572.5
$ws.Range("N11").Value = -5898850.5
$ws.Range("H23").Value = 3099.6667
$ws.Range("J23").Value = 3099.6667
$ws.Range("L23").Value = 3099.6667
$ws.Range("N23").Value = -3545.6667
$ws.Range("H70").Value = 7888.5
$ws.Range("I70").Value = 7337.5713
$ws.Range("K70").Value = 7337.5713
$ws.Range("M70").Value = -7067.5713
$ws.Range("H73").Value = 7888.5
$ws.Range("I73").Value = 7337.5713
$ws.Range("K73").Value = 7337.5713
$ws.Range("M73").Value = -6401.5713
$ws.Range("H95").Value = 100045640
$ws.Range("J95").Value = 100045640
$ws.Range("L95").Value = 100045640
$ws.Range("N95").Value = -100051132
$ws.Range("H102").Value = 2493.88
$ws.Range("I102").Value = 1687.4706
$ws.Range("K102").Value = 1687.4706
$ws.Range("M102").Value = -65.4706000000001
$ws.Range("H126").Value = 3326.6667
$ws.Range("J126").Value = 3961.7778
$ws.Range("L126").Value = 11885.3334
$ws.Range("N126").Value = -16825.3334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 25567.889
$ws.Range("J20").Value = 25567.889
$ws.Range("L20").Value = 25567.889
$ws.Range("N20").Value = -26019.889
$ws.Range("H46").Value = 4022.2122
$ws.Range("J46").Value = 8515.6
$ws.Range("L46").Value = 8515.6
$ws.Range("N46").Value = -8891.6
$ws.Range("H56").Value = 24834
$ws.Range("J56").Value = 44495
$ws.Range("L56").Value = 44495
$ws.Range("N56").Value = -45877
$ws.Range("H100").Value = 2927
$ws.Range("I100").Value = 1747.25
$ws.Range("K100").Value = 1747.25
$ws.Range("M100").Value = -1206.25
$ws.Range("H109").Value = 98274.664
$ws.Range("J109").Value = 98274.664
$ws.Range("L109").Value = 98274.664
$ws.Range("N109").Value = -101048.664
$ws.Range("H122").Value = 5043.0645
$ws.Range("I122").Value = 4737.2383
$ws.Range("K122").Value = 14211.7149
$ws.Range("M122").Value = -11761.7149
$ws.Range("H123").Value = 87980
$ws.Range("J123").Value = 87980
$ws.Range("L123").Value = 87980
$ws.Range("N123").Value = -97780
$ws.Range("H136").Value = 119316
$ws.Range("I136").Value = 18462.834
$ws.Range("J136").Value = 220169.17
$ws.Range("K136").Value = 55388.50199999999
$ws.Range("L136").Value = 660507.51
$ws.Range("M136").Value = -52838.50199999999
$ws.Range("N136").Value = -665607.51

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 104860
$ws.Range("J109").Value = 104860
$ws.Range("L109").Value = 104860
$ws.Range("N109").Value = -107634
$ws.Range("H113").Value = 474.6875
$ws.Range("J113").Value = 947.5
$ws.Range("L113").Value = 2842.5
$ws.Range("N113").Value = -7182.5
$ws.Range("H136").Value = 4448.4
$ws.Range("I136").Value = 4871.6
$ws.Range("J136").Value = 3178.8
$ws.Range("K136").Value = 14614.8
$ws.Range("L136").Value = 9536.400000000001
$ws.Range("M136").Value = -12064.8
$ws.Range("N136").Value = -14636.4
